$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (mean)
$ws.Range("B3").Value = 30629.02242373102
$ws.Range("D3").Value = 2001.12292375598
$ws.Range("E3").Value = 5766.766797521235

# Row 4 (std)
$ws.Range("B4").Value = 13511.2030964431
$ws.Range("D4").Value = 1247.804389540655
$ws.Range("E4").Value = 4278.85471380049

# Row 5 (min)
$ws.Range("B5").Value = 5285.039000000002
$ws.Range("D5").Value = 5.001
$ws.Range("E5").Value = 169.012

# Row 6 (25%)
$ws.Range("B6").Value = 19549.02975
$ws.Range("D6").Value = 814.003
$ws.Range("E6").Value = 2653

# Row 7 (50%)
$ws.Range("B7").Value = 27817.31150000003
$ws.Range("D7").Value = 1943
$ws.Range("E7").Value = 4615.0175

# Row 8 (75%)
$ws.Range("B8").Value = 41079.22675
$ws.Range("D8").Value = 3235.001
$ws.Range("E8").Value = 7265.011999999999

# Row 9 (max)
$ws.Range("B9").Value = 78071.52599999974
$ws.Range("D9").Value = 6008.009
$ws.Range("E9").Value = 35450.013

# Row 10 (Total)
$ws.Range("F10").Value = 16098614185.915

# Row 11 (Residential)
$ws.Range("G11").Value = 0.7463879318832842

# Row 12 (Community)
$ws.Range("F12").Value = 1051790208.726
$ws.Range("G12").Value = 0.06533420806160024

# Row 13 (IGA)
$ws.Range("F13").Value = 3031012628.777001
$ws.Range("G13").Value = 0.1882778600551155
